$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 71.125
$ws.Range("I5").Value = 71.125
$ws.Range("K5").Value = 71.125
$ws.Range("M5").Value = 43.875

$ws.Range("H40").Value = 5584
$ws.Range("J40").Value = 8392.866
$ws.Range("L40").Value = 8392.866
$ws.Range("N40").Value = -8742.866

$ws.Range("H93").Value = 40000
$ws.Range("J93").Value = 40000
$ws.Range("L93").Value = 40000
$ws.Range("N93").Value = -44992

$ws.Range("H101").Value = 734.25
$ws.Range("I101").Value = 474.75
$ws.Range("J101").Value = 993.75
$ws.Range("K101").Value = 1424.25
$ws.Range("L101").Value = 2981.25
$ws.Range("M101").Value = 197.75
$ws.Range("N101").Value = -6225.25

$ws.Range("H107").Value = 78643
$ws.Range("I107").Value = 85129.914
$ws.Range("K107").Value = 85129.914
$ws.Range("M107").Value = -83209.914

$ws.Range("H132").Value = 2485.818
$ws.Range("I132").Value = 1538.2222
$ws.Range("K132").Value = 4614.6666
$ws.Range("M132").Value = -2084.6666

$ws.Range("H138").Value = 3663.5312
$ws.Range("I138").Value = 1913.7727
$ws.Range("J138").Value = 4580.0713
$ws.Range("K138").Value = 5741.3181
$ws.Range("L138").Value = 13740.2139
$ws.Range("M138").Value = -601.3181000000004
$ws.Range("N138").Value = -24020.2139

$ws.Range("H141").Value = 6383.857
$ws.Range("J141").Value = 6383.857
$ws.Range("L141").Value = 19151.571
$ws.Range("N141").Value = -29511.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3713.5715
$ws.Range("I45").Value = 3759
$ws.Range("K45").Value = 3759
$ws.Range("M45").Value = -3382

$ws.Range("H102").Value = 1415.8276
$ws.Range("I102").Value = 1446.4286
$ws.Range("K102").Value = 1446.4286
$ws.Range("M102").Value = 175.5714

$ws.Range("H110").Value = 334627
$ws.Range("I110").Value = 358443.22
$ws.Range("K110").Value = 358443.22
$ws.Range("M110").Value = -356398.22

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2109.4707
$ws.Range("I99").Value = 1724.1333
$ws.Range("K99").Value = 1724.1333
$ws.Range("M99").Value = -226.1333

$ws.Range("H100").Value = 35000
$ws.Range("J100").Value = 35000
$ws.Range("L100").Value = 35000
$ws.Range("N100").Value = -37164

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4218.857
$ws.Range("I31").Value = 2114.625
$ws.Range("K31").Value = 2114.625
$ws.Range("M31").Value = -1819.625

$ws.Range("H34").Value = 4218.857
$ws.Range("I34").Value = 2114.625
$ws.Range("K34").Value = 2114.625
$ws.Range("M34").Value = -1912.625

$ws.Range("H55").Value = 16488
$ws.Range("I55").Value = 14290.667
$ws.Range("K55").Value = 14290.667
$ws.Range("M55").Value = -13975.667

$ws.Range("H62").Value = 3379.25
$ws.Range("I62").Value = 3309.889
$ws.Range("K62").Value = 3309.889
$ws.Range("M62").Value = -2685.889

$ws.Range("H65").Value = 3379.25
$ws.Range("I65").Value = 3309.889
$ws.Range("K65").Value = 16549.445
$ws.Range("M65").Value = -13429.445

$ws.Range("H94").Value = 4266
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 4266
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 4266
$ws.Range("M94").Value = ""
$ws.Range("N94").Value = -5168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 962.6875
$ws.Range("I5").Value = 758.6667
$ws.Range("K5").Value = 2276.0001
$ws.Range("M5").Value = -2164.0001

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""

$ws.Range("H135").Value = 962.6875
$ws.Range("I135").Value = 758.6667
$ws.Range("K135").Value = 6828.0003
$ws.Range("M135").Value = -4293.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 14400000
$ws.Range("J11").Value = 10000000
$ws.Range("L11").Value = 10000000
$ws.Range("N11").Value = -10000278

$ws.Range("H12").Value = 15000000
$ws.Range("I12").Value = 15000000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 15000000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -14999860
$ws.Range("N12").Value = ""

$ws.Range("H36").Value = 7563.2
$ws.Range("I36").Value = 4454
$ws.Range("K36").Value = 4454
$ws.Range("M36").Value = -3969

$ws.Range("H43").Value = 21749.625
$ws.Range("J43").Value = 30399.6
$ws.Range("L43").Value = 30399.6
$ws.Range("N43").Value = -30701.6

$ws.Range("H93").Value = 43999
$ws.Range("J93").Value = 43999
$ws.Range("L93").Value = 43999
$ws.Range("N93").Value = -47743

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5791.25
$ws.Range("I16").Value = 5221.8335
$ws.Range("J16").Value = 7499.5
$ws.Range("K16").Value = 5221.8335
$ws.Range("L16").Value = 7499.5
$ws.Range("M16").Value = -5051.8335
$ws.Range("N16").Value = -7839.5

$ws.Range("H20").Value = 10000000
$ws.Range("I20").Value = 10000000
$ws.Range("K20").Value = 10000000
$ws.Range("M20").Value = -9999774

$ws.Range("H22").Value = 854.8889
$ws.Range("J22").Value = 799.3333
$ws.Range("L22").Value = 799.3333
$ws.Range("N22").Value = -1389.3333

$ws.Range("H27").Value = 854.8889
$ws.Range("J27").Value = 799.3333
$ws.Range("L27").Value = 799.3333
$ws.Range("N27").Value = -1013.3333

$ws.Range("H93").Value = 50001296
$ws.Range("I93").Value = 55556940
$ws.Range("K93").Value = 55556940
$ws.Range("M93").Value = -55555692

$ws.Range("H132").Value = 4699.524
$ws.Range("I132").Value = 3650.2
$ws.Range("K132").Value = 10950.6
$ws.Range("M132").Value = -8420.599999999999

$ws.Range("H136").Value = 3554.8333
$ws.Range("J136").Value = 3984.9167
$ws.Range("L136").Value = 11954.7501
$ws.Range("N136").Value = -17054.7501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""

$ws.Range("H113").Value = 870.95654
$ws.Range("I113").Value = 932.5714
$ws.Range("K113").Value = 2797.7142
$ws.Range("M113").Value = -627.7142000000003

$ws.Range("H126").Value = 2158.889
$ws.Range("I126").Value = 2002.909
$ws.Range("J126").Value = 2587.8333
$ws.Range("K126").Value = 6008.727000000001
$ws.Range("L126").Value = 7763.499899999999
$ws.Range("M126").Value = -3538.727000000001
$ws.Range("N126").Value = -12703.4999

$ws.Range("H131").Value = 146821.33
$ws.Range("J131").Value = 146821.33
$ws.Range("L131").Value = 146821.33
$ws.Range("N131").Value = -156901.33
